# Apply the FHIR StructureDefinition publishing-run update:
#  - bump Version / Date on the "Metadata" sheet
#  - set Publisher, replace the duplicated "Contact" row with a
#    "Jurisdiction" row (the second duplicate row is removed, shifting
#    everything below it up by one row)
#  - refresh the root Extension's Short/Definition text on the
#    "Elements" sheet to the real title/description

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was previously blank
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 ("Contact" / "No display for ContactDetail") becomes the new
# "Jurisdiction" / "United States of America" row
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a second, duplicate "Contact" row - delete it so the rest of
# the table (Description, Purpose, Copyright, ...) shifts up by one row
$meta.Rows("11").Delete()

$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: give it the real Short/Definition text instead of
# the generic placeholders
$elements.Range("K2").Value = "Employee Year Of Service"
$elements.Range("L2").Value = "Net credited service in years for the employee"
